# CR_CronogramaResponsabilidades.xlsx edit script
# - Fills in the previously-empty F10:G10 dates
# - Inserts a new row at row 12 (Documento Esquema de Repositorios)
# - Fills the previously-empty row 11 with "Documento de Stakeholders" data
# - Fills the newly inserted row 12 with "Documento Esquema de Repositorios" data
# - Leaves the final selection on E11 (matches the author's last click before saving)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) F10:G10 were blank date cells; the author stamped them with the same
#    "last modified" date (2018-05-04 => serial 43224) used on the surrounding rows.
$ws.Range("F10").Value = 43224
$ws.Range("G10").Value = 43224

# 2) Insert a brand-new row above the old row 12 ("Caracteristicas - CAPITULO 1"),
#    pushing it (and everything below) down to row 13.
$ws.Rows.Item(12).Insert()

# Row 11 kept its original (blank-row) formatting, so just copy that over to the
# freshly inserted row 12 to match (fill color / borders / date number format).
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Fill row 11 - new task "Documento de Stakeholders"
$ws.Range("A11").Value = "Documento de Stakeholders"
$ws.Range("B11").Value = "SANTA CRUZ MIÑANO, ANGEL ANTONIO"
$ws.Range("C11").Value = "FINALIZADO"
$ws.Range("D11").Value = 'C:\Users\angel\Desktop\gasAPP\aplicacionesMOB_cc75\appGAS - Analisis\STH_StakeHolders'
$ws.Range("F11").Value = 43224
$ws.Range("G11").Value = 43224

# 4) Fill the newly inserted row 12 - new task "Documento Esquema de Repositorios"
$ws.Range("A12").Value = "Documento Esquema de Repositorios"
$ws.Range("B12").Value = "SANTA CRUZ MIÑANO, ANGEL ANTONIO"
$ws.Range("C12").Value = "FINALIZADO"
$ws.Range("D12").Value = 'C:\Users\angel\Desktop\gasAPP\aplicacionesMOB_cc75\appGAS - Gestion\PGC_PlanDeGestionDeConfiguracion\EI_EsquemaDeRepositorio'

# 5) Leave the selection where the author left it before saving.
$ws.Range("E11").Select()
